$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pedro Vyctor"
$ws.Range("B2").Value = "Carvalho"
$ws.Range("C2").Value = "pedro.vyctor00@gmail.com"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "Site perfeito!!"
